{"js": "// The document ends with two inline pictures, each in its own paragraph.\n// This edit:\n//   1) marks the first picture's run as \"do not spell/grammar check\"\n//      (adds <w:noProof/> to that run's rPr), and\n//   2) removes the second (last) picture entirely, along with the\n//      paragraph that holds it.\nconst body = context.document.body;\nconst pics = body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nif (pics.items.length >= 2) {\n  const firstPicture = pics.items[0];\n  const lastPicture = pics.items[pics.items.length - 1];\n\n  // 1) Flag the first picture's run as noProof.\n  const firstRange = firstPicture.getRange();\n  firstRange.hasNoProofing = true;\n\n  // 2) Delete the paragraph containing the last picture.\n  const lastParagraph = lastPicture.paragraph;\n  lastParagraph.delete();\n\n  await context.sync();\n}\n", "ps1": "# The document ends with two inline pictures, each in its own paragraph.\n# This edit:\n#   1) marks the first picture's run as \"do not spell/grammar check\"\n#      (adds <w:noProof/> to that run's rPr), and\n#   2) removes the second (last) picture entirely, along with the\n#      paragraph that holds it.\n$d = $word.ActiveDocument\n\nif ($d.InlineShapes.Count -ge 2) {\n    $firstShape = $d.InlineShapes.Item(1)\n    $lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)\n\n    # 1) Flag the first picture's run as noProof.\n    $firstShape.Range.NoProofing = 1\n\n    # 2) Delete the paragraph containing the last picture.\n    $lastParagraph = $lastShape.Range.Paragraphs.Item(1)\n    $lastParagraph.Range.Delete()\n}\n"}
